# Update the "Duration" column (E) values on the "Test status" and
# "Passed tests" sheets to reflect the new protractor run timings.

$wb = $excel.ActiveWorkbook

# row -> new value for column E (Duration)
$rowValues = @(
    @{ Row = 2; Value = 2.232 },
    @{ Row = 3; Value = 1.483 },
    @{ Row = 4; Value = 1.389 },
    @{ Row = 5; Value = 1.698 },
    @{ Row = 6; Value = 7.462 }
)

$sheetNames = @("Test status", "Passed tests")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $rowValues) {
        $ws.Cells.Item($entry.Row, 5).Value = $entry.Value
    }
}
